# StaticData TypeTest 추가 (#70)
# Adds a new "TypeTest" worksheet (after ClassListTest) containing
# min/max boundary values for a set of primitive .NET types.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip, make it active ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "TypeTest"

# --- note row ---
$ws.Range("A1").Value = "C10"

# --- side note near the TimeSpan column ---
$ws.Range("L9").Value = "ISO 8601 참고"
$ws.Range("M9").Value = "1일 10초"

# --- header row ---
$ws.Range("C10").Value = "Id"
$ws.Range("D10").Value = "UIntValue"
$ws.Range("E10").Value = "ShortValue"
$ws.Range("F10").Value = "LongValue"
$ws.Range("G10").Value = "FloatValue"
$ws.Range("H10").Value = "DoubleValue"
$ws.Range("I10").Value = "CharValue"
$ws.Range("J10").Value = "StringValue"
$ws.Range("K10").Value = "EnumValue"
$ws.Range("L10").Value = "DateTimeValue"
$ws.Range("M10").Value = "TimeSpanValue"

# --- row 11 : MinValue boundaries ---
$ws.Range("C11").Value = -2147483648.0
$ws.Range("D11").Value = 0.0
$ws.Range("E11").Value = -32768.0
$ws.Range("F11").Value = -9223372036854770000.0
$ws.Range("G11").Value = -340282350000000000000000000000000000000.0
$ws.Range("G11").NumberFormat = "0.00E+00"
$ws.Range("H11").Value = "-1.7976931348623157E+308"
$ws.Range("I11").Value = "0x00"
$ws.Range("J11").Value = "MinValue"
$ws.Range("K11").Value = "Sunday"
$ws.Range("L11").Value = "0001-01-01T00:00:00Z"
$ws.Range("M11").Value = "-10675199.02:48:05.4775808"

# --- row 12 : sample / typical values ---
$ws.Range("C12").Value = 1001.0
$ws.Range("D12").Value = 100.0
$ws.Range("E12").Value = 100.0
$ws.Range("F12").Value = 100.0
$ws.Range("G12").Value = 3.1415920000000002
$ws.Range("H12").Value = 3.1415926535896999
$ws.Range("I12").Value = "a"
$ws.Range("J12").Value = "ㅋㅋㅋ"
$ws.Range("K12").Value = "Monday"
$ws.Range("L12").Value = "1986-05-26T01:05:00+09:00"
$ws.Range("M12").Value = "1.00:00:10"

# --- row 13 : MaxValue boundaries ---
$ws.Range("C13").Value = 2147483647.0
$ws.Range("D13").Value = 4294967295.0
$ws.Range("E13").Value = 32767.0
$ws.Range("F13").Value = 9223372036854770000.0
$ws.Range("G13").Value = 340282350000000000000000000000000000000.0
$ws.Range("G13").NumberFormat = "0.00E+00"
$ws.Range("H13").Value = "1.7976931348623157E+308"
$ws.Range("I13").Value = "0xFFFF"
$ws.Range("J13").Value = "MaxValue"
$ws.Range("K13").Value = "Saturday"
$ws.Range("L13").Value = "9999-12-31T23:59:59Z"
$ws.Range("M13").Value = "10675199.02:48:05.4775807"

# --- layout: size columns to fit their (now entered) contents ---
$ws.Columns.Item(1).AutoFit()
$ws.Range("C1:M13").Columns.AutoFit()

# --- selection matches the author's last position on the new sheet ---
$ws.Range("E7").Select()
